$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly scoreboard rows (week of 2024-07-06, Excel serial date 45479)
# Columns: A Participant, B Date, C Workout Type, D Total Duration,
#          E Total Distance, F Total Elevation, G Zone1, H Zone2, I Zone3,
#          J Zone4, K Zone5, L Workout Level, M Week
$rows = @(
    @("Matt",     45479, "Run",     69, 5.03,  285, 9,  35, 8,  3,  0, "Sauntering Hippo", 4),
    @("Steven",   45479, "Workout", 40, 0,     0,   40, 0,  0,  0,  0, "Mighty Monkey",    4),
    @("Steven",   45479, "Walk",    15, 0.76,  16,  16, 0,  0,  0,  0, "Mighty Monkey",    4),
    @("Eric",     45479, "Run",     31, 3.04,  72,  0,  1,  10, 16, 0, "Agile Antelope",   4),
    @("Steven",   45479, "Walk",    17, 0.57,  89,  17, 0,  0,  0,  0, "Mighty Monkey",    4),
    @("Jeremiah", 45479, "Ride",    30, 11.24, 0,   0,  24, 5,  0,  0, "Agile Antelope",   4),
    @("Jeremiah", 45479, "Workout", 23, 0,     0,   21, 2,  0,  0,  0, "Agile Antelope",   4),
    @("Steven",   45479, "Walk",    44, 2.07,  26,  44, 0,  0,  0,  0, "Mighty Monkey",    4)
)

$startRow = 164
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r - 1, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
    $ws.Cells.Item($r, 13).Value = $data[12]
}

$ws.Range("A172").Select()
